$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCodes = @(
    "000004.SZ",
    "600814.SH",
    "600792.SH",
    "600794.SH",
    "600796.SH",
    "600797.SH",
    "600798.SH",
    "600800.SH",
    "600801.SH",
    "600802.SH",
    "600805.SH",
    "600808.SH",
    "600809.SH",
    "600811.SH",
    "600815.SH",
    "600789.SH",
    "600817.SH",
    "600818.SH",
    "600819.SH",
    "600822.SH",
    "600824.SH",
    "600825.SH",
    "600826.SH",
    "600830.SH",
    "600833.SH",
    "600834.SH",
    "600836.SH",
    "600838.SH",
    "600790.SH",
    "600787.SH",
    "600513.SH",
    "600765.SH",
    "600742.SH",
    "600745.SH",
    "600746.SH",
    "600749.SH",
    "600750.SH",
    "600751.SH",
    "600753.SH",
    "600756.SH",
    "600757.SH",
    "600761.SH",
    "600763.SH",
    "600764.SH",
    "600766.SH",
    "600784.SH",
    "600767.SH",
    "600768.SH",
    "600770.SH",
    "600771.SH"
)

$newNames = @(
    "ST国华",
    "杭州解百",
    "云煤能源",
    "保税科技",
    "钱江生化",
    "浙大网新",
    "宁波海运",
    "渤海化学",
    "华新水泥",
    "福建水泥",
    "悦达投资",
    "马钢股份",
    "山西汾酒",
    "东方集团",
    "厦工股份",
    "鲁抗医药",
    "宇通重工",
    "中路股份",
    "耀皮玻璃",
    "上海物贸",
    "益民集团",
    "新华传媒",
    "兰生股份",
    "香溢融通",
    "第一医药",
    "申通地铁",
    "上海易连",
    "上海九百",
    "轻纺城",
    "中储股份",
    "联环药业",
    "中航重机",
    "一汽富维",
    "闻泰科技",
    "江苏索普",
    "西藏旅游",
    "江中药业",
    "海航科技",
    "东方银星",
    "浪潮软件",
    "长江传媒",
    "安徽合力",
    "通策医疗",
    "中国海防",
    "*ST园城",
    "鲁银投资",
    "*ST运盛",
    "宁波富邦",
    "综艺股份",
    "广誉远"
)

$newF = @(
    8.83,
    7.1,
    3.95,
    3.28,
    5.21,
    5.41,
    4.19,
    3.95,
    21.47,
    6.39,
    4.83,
    3.9,
    254.95,
    2.89,
    2.9,
    7.04,
    7.95,
    17.51,
    4.84,
    8.57,
    4.16,
    3.8,
    7.55,
    5.29,
    9.82,
    8.48,
    5.79,
    10.1,
    4.08,
    5.72,
    9.949999999999999,
    40.03,
    8.6,
    67.54000000000001,
    10.7,
    12.18,
    13.96,
    2.73,
    11.67,
    12.36,
    5.38,
    9.869999999999999,
    125.88,
    22.37,
    6,
    6.15,
    4.02,
    8.16,
    10.03,
    29.4
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]
    $ws.Cells.Item($row, 2).Value = $newNames[$i]
    $ws.Cells.Item($row, 6).Value = $newF[$i]
}
